# This workbook holds one weekly price-report row per line. A new weekly
# report (row 39's new values below) is inserted at the top of the data,
# which pushes every existing row's report down by one week: row N's old
# data becomes row N+1's data, for every row from the (old) last row back
# up through row 40. The row that falls off the bottom (old row 156) is
# appended as the new last row, 157.
#
# Columns A, B, C, E, F, G, H, N, Q, R are constant for every data row in
# this sheet (same market/region/product/unit/etc.), so only D, I, J, K,
# L, M, O, P actually carry per-week data and need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataCols = @(4, 9, 10, 11, 12, 13, 15, 16)  # D, I, J, K, L, M, O, P
$firstDataRow = 40
$lastDataRow = 156
$newLastRow = 157

# Capture the old last row's values before anything is overwritten -
# these become the new appended row.
$carryOver = @{}
foreach ($c in $dataCols) {
    $carryOver[$c] = $ws.Cells.Item($lastDataRow, $c).Value2
}

# Shift every row's data down by one, walking from the bottom up so each
# source row is read before it gets overwritten.
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    foreach ($c in $dataCols) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($r - 1, $c).Value2
    }
}

# Row 39 gets the brand-new weekly report.
$ws.Range("D39").Value = 44414
$ws.Range("J39").Value = 250
$ws.Range("K39").Value = 3500
$ws.Range("L39").Value = 4000
$ws.Range("M39").Value = 3760
$ws.Range("P39").Value = 1253

# Append the row that was pushed off the bottom as new row 157, filling
# in the constant columns too since this is a brand-new row.
$ws.Range("A$newLastRow").Value = 3
$ws.Range("B$newLastRow").Value = "Femacal de La Calera"
$ws.Range("C$newLastRow").Value = "Coquimbo"
$ws.Range("D$newLastRow").Value = $carryOver[4]
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastDataRow").NumberFormat
$ws.Range("E$newLastRow").Value = 5
$ws.Range("F$newLastRow").Value = 100112012
$ws.Range("G$newLastRow").Value = "Espinaca"
$ws.Range("H$newLastRow").Value = "Sin especificar"
$ws.Range("I$newLastRow").Value = $carryOver[9]
$ws.Range("J$newLastRow").Value = $carryOver[10]
$ws.Range("K$newLastRow").Value = $carryOver[11]
$ws.Range("L$newLastRow").Value = $carryOver[12]
$ws.Range("M$newLastRow").Value = $carryOver[13]
$ws.Range("N$newLastRow").Value = "$/docena de atados (3 kilos)"
$ws.Range("O$newLastRow").Value = $carryOver[15]
$ws.Range("P$newLastRow").Value = $carryOver[16]
$ws.Range("Q$newLastRow").Value = 3
$ws.Range("R$newLastRow").Value = "Hortaliza"
